# Append 10 coded-segment rows (sheet rows 284-293) that mirror the
# existing row-283 layout/styling exactly, per "run latest mex files.
# update qa checks on segments, drugs, bacteria."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- row 284 ----
$ws.Range("A284").Value2 = "'●"
$ws.Range("B284").Value2 = "'"
$ws.Range("C284").Value2 = "'"
$ws.Range("D284").Value2 = "'2863"
$ws.Range("E284").Value2 = "'Event year"
$ws.Range("F284").Value2 = "'1: 1933"
$ws.Range("G284").Value2 = "'1: 1936"
$ws.Range("I284").Value2 = "'2000"
$ws.Range("L284").Value2 = "'Sonia"
$ws.Range("M284").Value2 = "'11/14/18 13:12:00"
$ws.Range("H284").Value2 = 0
$ws.Range("J284").Value2 = 4
$ws.Range("K284").Value2 = 0.011254
$ws.Range("A283:M283").Copy()
$ws.Range("A284:M284").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(284).RowHeight = 16

# ---- row 285 ----
$ws.Range("A285").Value2 = "'●"
$ws.Range("B285").Value2 = "'"
$ws.Range("C285").Value2 = "'"
$ws.Range("D285").Value2 = "'2863"
$ws.Range("E285").Value2 = "'Event year"
$ws.Range("F285").Value2 = "'1: 1942"
$ws.Range("G285").Value2 = "'1: 1945"
$ws.Range("I285").Value2 = "'2013"
$ws.Range("L285").Value2 = "'Sonia"
$ws.Range("M285").Value2 = "'11/14/18 13:12:00"
$ws.Range("H285").Value2 = 0
$ws.Range("J285").Value2 = 4
$ws.Range("K285").Value2 = 0.011254
$ws.Range("A283:M283").Copy()
$ws.Range("A285:M285").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(285).RowHeight = 16

# ---- row 286 ----
$ws.Range("A286").Value2 = "'●"
$ws.Range("B286").Value2 = "'"
$ws.Range("C286").Value2 = "'"
$ws.Range("D286").Value2 = "'5362"
$ws.Range("E286").Value2 = "'Event month"
$ws.Range("F286").Value2 = "'1: 1890"
$ws.Range("G286").Value2 = "'1: 1893"
$ws.Range("I286").Value2 = "'June"
$ws.Range("L286").Value2 = "'Sonia"
$ws.Range("M286").Value2 = "'11/14/18 13:17:00"
$ws.Range("H286").Value2 = 0
$ws.Range("J286").Value2 = 4
$ws.Range("K286").Value2 = 0.011473
$ws.Range("A283:M283").Copy()
$ws.Range("A286:M286").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(286).RowHeight = 16

# ---- row 287 ----
$ws.Range("A287").Value2 = "'●"
$ws.Range("B287").Value2 = "'"
$ws.Range("C287").Value2 = "'"
$ws.Range("D287").Value2 = "'5362"
$ws.Range("E287").Value2 = "'Event month"
$ws.Range("F287").Value2 = "'1: 1903"
$ws.Range("G287").Value2 = "'1: 1909"
$ws.Range("I287").Value2 = "'January"
$ws.Range("L287").Value2 = "'Sonia"
$ws.Range("M287").Value2 = "'11/14/18 13:17:00"
$ws.Range("H287").Value2 = 0
$ws.Range("J287").Value2 = 7
$ws.Range("K287").Value2 = 0.020077
$ws.Range("A283:M283").Copy()
$ws.Range("A287:M287").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(287).RowHeight = 16

# ---- row 288 ----
$ws.Range("A288").Value2 = "'●"
$ws.Range("B288").Value2 = "'"
$ws.Range("C288").Value2 = "'"
$ws.Range("D288").Value2 = "'5362"
$ws.Range("E288").Value2 = "'Event year"
$ws.Range("F288").Value2 = "'1: 1895"
$ws.Range("G288").Value2 = "'1: 1898"
$ws.Range("I288").Value2 = "'2009"
$ws.Range("L288").Value2 = "'Sonia"
$ws.Range("M288").Value2 = "'11/14/18 13:17:00"
$ws.Range("H288").Value2 = 0
$ws.Range("J288").Value2 = 4
$ws.Range("K288").Value2 = 0.011473
$ws.Range("A283:M283").Copy()
$ws.Range("A288:M288").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(288).RowHeight = 16

# ---- row 289 ----
$ws.Range("A289").Value2 = "'●"
$ws.Range("B289").Value2 = "'"
$ws.Range("C289").Value2 = "'"
$ws.Range("D289").Value2 = "'5362"
$ws.Range("E289").Value2 = "'Event year"
$ws.Range("F289").Value2 = "'1: 1911"
$ws.Range("G289").Value2 = "'1: 1914"
$ws.Range("I289").Value2 = "'2010"
$ws.Range("L289").Value2 = "'Sonia"
$ws.Range("M289").Value2 = "'11/14/18 13:18:00"
$ws.Range("H289").Value2 = 0
$ws.Range("J289").Value2 = 4
$ws.Range("K289").Value2 = 0.011473
$ws.Range("A283:M283").Copy()
$ws.Range("A289:M289").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(289).RowHeight = 16

# ---- row 290 ----
$ws.Range("A290").Value2 = "'●"
$ws.Range("B290").Value2 = "'"
$ws.Range("C290").Value2 = "'"
$ws.Range("D290").Value2 = "'5362"
$ws.Range("E290").Value2 = "'B"
$ws.Range("F290").Value2 = "'1: 1911"
$ws.Range("G290").Value2 = "'1: 1914"
$ws.Range("I290").Value2 = "'2010"
$ws.Range("L290").Value2 = "'Sonia"
$ws.Range("M290").Value2 = "'11/14/18 13:18:00"
$ws.Range("H290").Value2 = 0
$ws.Range("J290").Value2 = 4
$ws.Range("K290").Value2 = 0.011473
$ws.Range("A283:M283").Copy()
$ws.Range("A290:M290").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(290).RowHeight = 16

# ---- row 291 ----
$ws.Range("A291").Value2 = "'●"
$ws.Range("B291").Value2 = "'"
$ws.Range("C291").Value2 = "'"
$ws.Range("D291").Value2 = "'5362"
$ws.Range("E291").Value2 = "'B"
$ws.Range("F291").Value2 = "'1: 1903"
$ws.Range("G291").Value2 = "'1: 1909"
$ws.Range("I291").Value2 = "'January"
$ws.Range("L291").Value2 = "'Sonia"
$ws.Range("M291").Value2 = "'11/14/18 13:18:00"
$ws.Range("H291").Value2 = 0
$ws.Range("J291").Value2 = 7
$ws.Range("K291").Value2 = 0.020077
$ws.Range("A283:M283").Copy()
$ws.Range("A291:M291").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(291).RowHeight = 16

# ---- row 292 ----
$ws.Range("A292").Value2 = "'●"
$ws.Range("B292").Value2 = "'"
$ws.Range("C292").Value2 = "'"
$ws.Range("D292").Value2 = "'5362"
$ws.Range("E292").Value2 = "'A"
$ws.Range("F292").Value2 = "'1: 1895"
$ws.Range("G292").Value2 = "'1: 1898"
$ws.Range("I292").Value2 = "'2009"
$ws.Range("L292").Value2 = "'Sonia"
$ws.Range("M292").Value2 = "'11/14/18 13:18:00"
$ws.Range("H292").Value2 = 0
$ws.Range("J292").Value2 = 4
$ws.Range("K292").Value2 = 0.011473
$ws.Range("A283:M283").Copy()
$ws.Range("A292:M292").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(292).RowHeight = 16

# ---- row 293 ----
$ws.Range("A293").Value2 = "'●"
$ws.Range("B293").Value2 = "'"
$ws.Range("C293").Value2 = "'"
$ws.Range("D293").Value2 = "'5362"
$ws.Range("E293").Value2 = "'A"
$ws.Range("F293").Value2 = "'1: 1890"
$ws.Range("G293").Value2 = "'1: 1893"
$ws.Range("I293").Value2 = "'June"
$ws.Range("L293").Value2 = "'Sonia"
$ws.Range("M293").Value2 = "'11/14/18 13:18:00"
$ws.Range("H293").Value2 = 0
$ws.Range("J293").Value2 = 4
$ws.Range("K293").Value2 = 0.011473
$ws.Range("A283:M283").Copy()
$ws.Range("A293:M293").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(293).RowHeight = 16

